# Auto-generated: update market-data columns (H:N) across multiple sheets
# per scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 539.8
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 539.8
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 539.8
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -677.8
$ws.Range("H69").Value = 3716.4092
$ws.Range("I69").Value = 10001
$ws.Range("J69").Value = 3417.1428
$ws.Range("K69").Value = 30003
$ws.Range("L69").Value = 10251.4284
$ws.Range("M69").Value = -29129
$ws.Range("N69").Value = -11999.4284
$ws.Range("H72").Value = 3716.4092
$ws.Range("I72").Value = 10001
$ws.Range("J72").Value = 3417.1428
$ws.Range("K72").Value = 90009
$ws.Range("L72").Value = 30754.2852
$ws.Range("M72").Value = -85641
$ws.Range("N72").Value = -39490.2852
$ws.Range("H80").Value = 2699.3076
$ws.Range("I80").Value = 1393.3334
$ws.Range("J80").Value = 3818.7144
$ws.Range("K80").Value = 4180.0002
$ws.Range("L80").Value = 11456.1432
$ws.Range("M80").Value = -3182.0002
$ws.Range("N80").Value = -13452.1432
$ws.Range("H83").Value = 2699.3076
$ws.Range("I83").Value = 1393.3334
$ws.Range("J83").Value = 3818.7144
$ws.Range("K83").Value = 12540.0006
$ws.Range("L83").Value = 34368.4296
$ws.Range("M83").Value = -7548.000599999999
$ws.Range("N83").Value = -44352.4296
$ws.Range("H131").Value = 847.5
$ws.Range("J131").Value = 600
$ws.Range("L131").Value = 1800
$ws.Range("N131").Value = -11880

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1167.3334
$ws.Range("I5").Value = 250
$ws.Range("J5").Value = 1626
$ws.Range("K5").Value = 250
$ws.Range("L5").Value = 1626
$ws.Range("M5").Value = -138
$ws.Range("N5").Value = -1850
$ws.Range("H11").Value = 3666.6667
$ws.Range("I11").Value = 1000
$ws.Range("K11").Value = 1000
$ws.Range("M11").Value = -856
$ws.Range("H31").Value = 6540.3335
$ws.Range("I31").Value = 6540.3335
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6540.3335
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -6246.3335
$ws.Range("N31").ClearContents()
$ws.Range("H33").Value = 5999.5
$ws.Range("I33").Value = 5999.5
$ws.Range("K33").Value = 5999.5
$ws.Range("M33").Value = -5670.5
$ws.Range("H36").Value = 7777
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 7777
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 7777
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -8469
$ws.Range("H122").Value = 1884.1321
$ws.Range("I122").Value = 1713.8864
$ws.Range("K122").Value = 5141.6592
$ws.Range("M122").Value = -2691.6592

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1167.3334
$ws.Range("I4").Value = 250
$ws.Range("J4").Value = 1626
$ws.Range("K4").Value = 250
$ws.Range("L4").Value = 1626
$ws.Range("M4").Value = -135
$ws.Range("N4").Value = -1856
$ws.Range("H7").Value = 2750.75
$ws.Range("I7").Value = 1501.5
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 1501.5
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -1388.5
$ws.Range("N7").Value = -4226
$ws.Range("H75").Value = 8066.8184
$ws.Range("I75").Value = 4303.8887
$ws.Range("K75").Value = 4303.8887
$ws.Range("M75").Value = -3367.8887
$ws.Range("H78").Value = 8066.8184
$ws.Range("I78").Value = 4303.8887
$ws.Range("K78").Value = 12911.6661
$ws.Range("M78").Value = -8231.666100000002
$ws.Range("H102").Value = 14511.2
$ws.Range("I102").Value = 9278
$ws.Range("J102").Value = 18000
$ws.Range("K102").Value = 9278
$ws.Range("L102").Value = 18000
$ws.Range("M102").Value = -6033
$ws.Range("N102").Value = -24490

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1066.1666
$ws.Range("I2").Value = 679.4
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 679.4
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -566.4
$ws.Range("N2").Value = -3226
$ws.Range("H7").Value = 132
$ws.Range("I7").Value = 20
$ws.Range("J7").Value = 199.2
$ws.Range("K7").Value = 20
$ws.Range("L7").Value = 199.2
$ws.Range("M7").Value = 93
$ws.Range("N7").Value = -425.2
$ws.Range("H22").Value = 280.1111
$ws.Range("I22").Value = 295.85715
$ws.Range("J22").Value = 225
$ws.Range("K22").Value = 295.85715
$ws.Range("L22").Value = 225
$ws.Range("M22").Value = 54.14285000000001
$ws.Range("N22").Value = -925
$ws.Range("H108").Value = 39833.168
$ws.Range("J108").Value = 39833.168
$ws.Range("L108").Value = 39833.168
$ws.Range("N108").Value = -47513.168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 411.07693
$ws.Range("I8").Value = 411.07693
$ws.Range("K8").Value = 1233.23079
$ws.Range("M8").Value = -1094.23079
$ws.Range("H114").Value = 631.2857
$ws.Range("I114").Value = 319.45456
$ws.Range("J114").Value = 833.05884
$ws.Range("K114").Value = 958.36368
$ws.Range("L114").Value = 2499.17652
$ws.Range("M114").Value = 2295.63632
$ws.Range("N114").Value = -9007.176520000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1450184.1
$ws.Range("I22").Value = 16666916
$ws.Range("J22").Value = 971.5238000000001
$ws.Range("K22").Value = 16666916
$ws.Range("L22").Value = 971.5238000000001
$ws.Range("M22").Value = -16666621
$ws.Range("N22").Value = -1561.5238
$ws.Range("H27").Value = 1450184.1
$ws.Range("I27").Value = 16666916
$ws.Range("J27").Value = 971.5238000000001
$ws.Range("K27").Value = 16666916
$ws.Range("L27").Value = 971.5238000000001
$ws.Range("M27").Value = -16666809
$ws.Range("N27").Value = -1185.5238

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7033.5713
$ws.Range("I62").Value = 6500
$ws.Range("J62").Value = 7122.5
$ws.Range("K62").Value = 6500
$ws.Range("L62").Value = 7122.5
$ws.Range("M62").Value = -5876
$ws.Range("N62").Value = -8370.5
$ws.Range("H65").Value = 7033.5713
$ws.Range("I65").Value = 6500
$ws.Range("J65").Value = 7122.5
$ws.Range("K65").Value = 32500
$ws.Range("L65").Value = 35612.5
$ws.Range("M65").Value = -29380
$ws.Range("N65").Value = -41852.5
